$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E9").Value = 113.8
$ws.Range("E13").Value = 89.5
$ws.Range("F13").Value = 79.317
$ws.Range("P13").Value = 1022.478
$ws.Range("S13").Value = 4675
$ws.Range("Y13").Value = 272.396
